$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# remain text (matching the source data, which stores everything as text)
# otherwise Excel auto-converts them to numeric cells.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '60.380.43'
$ws.Range('E2').Value = '  +4.18%  '
$ws.Range('D3').Value = '2.432.53'
$ws.Range('E3').Value = '  +3.17%  '
$ws.Range('E4').Value = '  -0.01%  '
Set-TextValue 'D5' '556.70'
$ws.Range('E5').Value = '  +2.49%  '
Set-TextValue 'D6' '139.42'
$ws.Range('E6').Value = '  +3.58%  '
$ws.Range('E7').Value = '  -0.05%  '
Set-TextValue 'D8' '0.578'
$ws.Range('E8').Value = '  +3.14%  '
$ws.Range('E9').Value = '  +4.81%  '
Set-TextValue 'D10' '5.76'
$ws.Range('E10').Value = '  +3.99%  '
Set-TextValue 'D11' '0.359'
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('E12').Value = '  -2.18%  '
$ws.Range('E13').Value = '  +5.17%  '
$ws.Range('D14').Value = '2.865.44'
$ws.Range('E14').Value = '  +3.14%  '
$ws.Range('D15').Value = '60.294.79'
$ws.Range('E15').Value = '  +4.10%  '
$ws.Range('E16').Value = '  +4.18%  '
$ws.Range('D17').Value = '2.403.83'
$ws.Range('E17').Value = '  +1.22%  '
Set-TextValue 'D18' '11.41'
Set-TextValue 'D19' '4.42'
$ws.Range('E19').Value = '  +3.03%  '
Set-TextValue 'D20' '334.53'
$ws.Range('E20').Value = '  +1.20%  '
Set-TextValue 'D21' '6.78'
$ws.Range('E21').Value = '  +0.93%  '
Set-TextValue 'D22' '1.00'
$ws.Range('E22').Value = '  -0.06%  '
Set-TextValue 'D23' '65.43'
$ws.Range('E23').Value = '  +4.26%  '
$ws.Range('E24').Value = '  +3.42%  '
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +0.26%  '
$ws.Range('D28').Value = '0.0₃0791'
$ws.Range('E28').Value = '  +6.69%  '
$ws.Range('E29').Value = '  +2.24%  '
Set-TextValue 'D30' '6.33'
$ws.Range('E30').Value = '  +2.99%  '
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('B32').Value = 'SuiNetwork'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue 'D32' '1.04'
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D33' '18.78'
$ws.Range('E33').Value = '  +1.98%  '
$ws.Range('E35').Value = '  +6.25%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  +0.50%  '
Set-TextValue 'D39' '39.86'
$ws.Range('E39').Value = '  +1.09%  '
Set-TextValue 'D40' '0.421'
$ws.Range('E40').Value = '  +10.96%  '
Set-TextValue 'D41' '321.71'
$ws.Range('E41').Value = '  +11.11%  '
$ws.Range('E42').Value = '  +1.46%  '
Set-TextValue 'D43' '141.48'
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D44' '0.0527'
$ws.Range('E44').Value = '  +3.60%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D45' '0.0961'
$ws.Range('E45').Value = '  +1.71%  '
Set-TextValue 'D46' '19.65'
$ws.Range('E46').Value = '  +2.62%  '
Set-TextValue 'D47' '0.413'
$ws.Range('E47').Value = '  +8.17%  '
Set-TextValue 'D48' '0.574'
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('E49').Value = '  +1.85%  '
Set-TextValue 'D50' '17.96'
$ws.Range('E50').Value = '  +2.63%  '
$ws.Range('E51').Value = '  -0.12%  '
